$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 22:01"

# Country ranking changed: update names + stats for the affected rows
# (countries that swapped relative rank, plus stat refreshes for others)

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 2374777
$ws.Range("C4").Value = 18120
$ws.Range("D4").Value = 984277
$ws.Range("E4").Value = 1268005
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 248
$ws.Range("H4").Value = 122495

# Row 7: India
$ws.Range("A7").Value = "India"
$ws.Range("B7").Value = 440450
$ws.Range("C7").Value = 13540
$ws.Range("D7").Value = 248137
$ws.Range("E7").Value = 178298
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 312
$ws.Range("H7").Value = 14015

# Row 19: Francia
$ws.Range("A19").Value = "Francia"
$ws.Range("B19").Value = 160750
$ws.Range("C19").Value = 373
$ws.Range("D19").Value = 74372
$ws.Range("E19").Value = 56715
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 29663

# Row 21: Canada
$ws.Range("A21").Value = "Canada"
$ws.Range("B21").Value = 101568
$ws.Range("C21").Value = 231
$ws.Range("D21").Value = 64284
$ws.Range("E21").Value = 28850
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 8434

# Row 28: Egipto
$ws.Range("A28").Value = "Egipto"
$ws.Range("B28").Value = 56809
$ws.Range("C28").Value = 1576
$ws.Range("D28").Value = 15133
$ws.Range("E28").Value = 39398
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 85
$ws.Range("H28").Value = 2278

# Row 29: Suecia
$ws.Range("A29").Value = "Suecia"
$ws.Range("B29").Value = 56043
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 5053

# Row 96: Mauritania
$ws.Range("A96").Value = "Mauritania"
$ws.Range("B96").Value = 3121
$ws.Range("C96").Value = 137
$ws.Range("D96").Value = 905
$ws.Range("E96").Value = 2104
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 112

# Row 97: Republica de Africa Central
$ws.Range("A97").Value = "Republica de Africa Central"
$ws.Range("B97").Value = 2963
$ws.Range("C97").Value = 155
$ws.Range("D97").Value = 495
$ws.Range("E97").Value = 2438
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 7
$ws.Range("H97").Value = 30

# Row 98: Somalia
$ws.Range("A98").Value = "Somalia"
$ws.Range("B98").Value = 2812
$ws.Range("C98").Value = 33
$ws.Range("D98").Value = 818
$ws.Range("E98").Value = 1904
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 90

# Row 99: Guayana Francesa
$ws.Range("A99").Value = "Guayana Francesa"
$ws.Range("B99").Value = 2458
$ws.Range("C99").Value = 17
$ws.Range("D99").Value = 972
$ws.Range("E99").Value = 1478
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 8

# Row 103: Costa Rica
$ws.Range("A103").Value = "Costa Rica"
$ws.Range("B103").Value = 2277
$ws.Range("C103").Value = 64
$ws.Range("D103").Value = 1043
$ws.Range("E103").Value = 1222
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 12

# Row 104: Maldivas
$ws.Range("A104").Value = "Maldivas"
$ws.Range("B104").Value = 2217
$ws.Range("C104").Value = 14
$ws.Range("D104").Value = 1813
$ws.Range("E104").Value = 396
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 8

# Row 117: Guinea-Bisau
$ws.Range("A117").Value = "Guinea-Bisau"
$ws.Range("B117").Value = 1556
$ws.Range("C117").Value = 15
$ws.Range("D117").Value = 191
$ws.Range("E117").Value = 1346
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 19

# Row 129: Estado de Palestina
$ws.Range("A129").Value = "Estado de Palestina"
$ws.Range("B129").Value = 1001
$ws.Range("C129").Value = 168
$ws.Range("D129").Value = 442
$ws.Range("E129").Value = 556
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 3

# Row 153: Reunion
$ws.Range("A153").Value = "Reunion"
$ws.Range("B153").Value = 507
$ws.Range("C153").Value = 1
$ws.Range("D153").Value = 460
$ws.Range("E153").Value = 46
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 1

# Row 179: Monaco
$ws.Range("A179").Value = "Monaco"
$ws.Range("B179").Value = 101
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 94
$ws.Range("E179").Value = 3
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 4

# Row 180: Aruba
$ws.Range("A180").Value = "Aruba"
$ws.Range("B180").Value = 101
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 98
$ws.Range("E180").Value = 0
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 4

# Row 202: Fiyi
$ws.Range("A202").Value = "Fiyi"
$ws.Range("B202").Value = 18
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 18
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

# Row 203: Dominica
$ws.Range("A203").Value = "Dominica"
$ws.Range("B203").Value = 18
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 18
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0
